# Split the single run that holds the "Цель проекта" sentence into three
# runs so that the word "мобильной" is replaced by "компьютерной", with the
# newly typed word carrying its own (slightly different) run formatting
# (an explicit complex-script size that the surrounding, untouched runs do
# not have) -- exactly mirroring what Word does when you select a word
# mid-sentence and retype it.

$d = $word.ActiveDocument

$oldSentence = "Целью данного проекта является разработка мобильной игры, которая будет интересной и увлекательной для игр."

$findRange = $d.Content
$found = $findRange.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $findRange.Start
    $end = $findRange.End
    $target = $d.Range($start, $end)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r w:rsidRPr="0002276A">' +
        '<w:rPr><w:bCs/><w:sz w:val="28"/></w:rPr>' +
        '<w:t xml:space="preserve">Целью данного проекта является разработка </w:t>' +
        '</w:r>' +
        '<w:r>' +
        '<w:rPr><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
        '<w:t xml:space="preserve">компьютерной </w:t>' +
        '</w:r>' +
        '<w:r>' +
        '<w:rPr><w:bCs/><w:sz w:val="28"/></w:rPr>' +
        '<w:t>игры, которая будет интересной и увлекательной для игр.</w:t>' +
        '</w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'

    $target.InsertXML($xml)
} else {
    Write-Output "Target sentence not found"
}
